# Generate Report for Archive
$wb = $excel.ActiveWorkbook

# 1. Update status text "Ready for handoff" -> "In Translation"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C3").Value = "In Translation"

# 2. Narrow the "Status" columns (target stored width ~= 13.4101845877511
#    characters; the COM ColumnWidth setter quantizes to whole pixels, so
#    12.5 is the input that lands on the closest achievable stored width)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
